# Automatische test-sync: 2025-07-29 22:00:50
# Append Testmail #15 ("Leg dit even neer bij Koen.") to the Logs sheet,
# extend the conditional-formatting ranges to include the new row, and
# bump the Dashboard summary count for the matching category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 17

$ws.Cells.Item($newRow, 1).Value = "Leg dit even neer bij Koen."
$ws.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value = "Testmail #15: Leg dit even neer bij Koen."
$ws.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
# Column E (Antwoord) intentionally left blank - no automated reply for this mail.
$ws.Cells.Item($newRow, 6).Value = "2025-07-29 21:59:52"
$ws.Cells.Item($newRow, 7).Value = "Nee"
$ws.Cells.Item($newRow, 8).Value = "Ja"
$ws.Cells.Item($newRow, 9).Value = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J) from row 16 to row 17.
$cfColumns = "D", "G", "H", "I", "J"
foreach ($col in $cfColumns) {
    $oldRange = $ws.Range($col + "2:" + $col + "16")
    $newRange = $ws.Range($col + "2:" + $col + "17")
    $fcs = $oldRange.FormatConditions
    $fcCount = $fcs.Count()
    for ($i = 1; $i -le $fcCount; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for "Intern verzoek / Actie voor medewerker".
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 2).Value = 4
